$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$wb.Sheets.Item(1).Name = "Through 2022-10-14"

$ws.Range("A11").Value = "October (through 10-14)"

$ws.Range("B11").Value = 13
$ws.Range("C11").Value = 22
$ws.Range("D11").Value = 23
$ws.Range("E11").Value = 36
$ws.Range("F11").Value = 19
$ws.Range("G11").Value = 65
$ws.Range("H11").Value = 88
$ws.Range("I11").Value = 45

$ws.Range("B12").Value = 239
$ws.Range("C12").Value = 451
$ws.Range("D12").Value = 650
$ws.Range("E12").Value = 584
$ws.Range("F12").Value = 441
$ws.Range("G12").Value = 966
$ws.Range("H12").Value = 1335
$ws.Range("I12").Value = 1323
